$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "quantity" in column I, row 1
$ws.Range("I1").Value = "quantity"

# Add quantity value in I5
$ws.Range("I5").Value = 4

# Add quantity value in H6
$ws.Range("H6").Value = 3
$ws.Rows.Item(6).RowHeight = 13.8

# Update selection to match target
$ws.Range("G15").Select()
